$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (TPM-updated values). Target clusters "Neutrophils" and
# "Resolving-Mac" are replaced throughout by "ECs" / "MuSCs" target rows,
# adding one extra data row (7 data rows total instead of 6).
$data = @(
  @("ECs",  "Efna4", "Epha5", "ECs",   2, 1, 0.595027,             1.190054,  0.3721518773584837, 0.3173844652954502, 1, 0.5, 0.006466, 0.012932, 0.008493427970384656, 0.008493427970384656, 0.003847444582,      0.015389778328,  0.003160845164387706, 0.002695682094905955),
  @("ECs",  "Efna4", "Epha5", "MuSCs", 2, 1, 0.595027,             1.190054,  0.3721518773584837, 0.3173844652954502, 2, 1,   0.7548284999999999, 1.509657, 0.9915065720296153,  0.9915065720296153,  0.4491433378694999, 1.796573351478, 0.368991032194096,   0.3146887832005443),
  @("FAPs", "Efna4", "Epha5", "ECs",   3, 1, 0.5518016666666666,  1.655405,  0.3451171563299485, 0.441492428723751,  1, 0.5, 0.006466, 0.012932, 0.008493427970384656, 0.008493427970384656, 0.003567949576666667,0.02140769746,  0.002931227708632398, 0.003749784142835361),
  @("FAPs", "Efna4", "Epha5", "MuSCs", 3, 1, 0.5518016666666666,  1.655405,  0.3451171563299485, 0.441492428723751,  2, 1,   0.7548284999999999, 1.509657, 0.9915065720296153,  0.9915065720296153,  0.4165156243474999, 2.499093746085, 0.3421859286213161,  0.4377426445809156),
  @("MuSCs","Efna4", "Epha5", "ECs",   2, 1, 0.4520535,           0.904107,  0.2827309663115679, 0.2411231059807989, 1, 0.5, 0.006466, 0.012932, 0.008493427970384656, 0.008493427970384656, 0.002922977931,      0.011691911724,  0.002401355097364553, 0.002047961732643341),
  @("MuSCs","Efna4", "Epha5", "MuSCs", 2, 1, 0.4520535,           0.904107,  0.2827309663115679, 0.2411231059807989, 2, 1,   0.7548284999999999, 1.509657, 0.9915065720296153,  0.9915065720296153,  0.3412228653247499, 1.364891461299, 0.2803296112142034,  0.2390751442481555)
)

$numRows = $data.Count
$numCols = $data[0].Count

# Clear the old data region first (old sheet used rows 2-6, new sheet uses
# rows 2-7), then write the full replacement block in one shot.
$ws.Range("A2:T20").ClearContents()

$arr = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$rangeAddress = "A" + $startRow + ":T" + $endRow
$ws.Range($rangeAddress).Value2 = $arr
